# Auto-generated edit script applying scheduled market-price/profit data refresh
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
$wsALC.Range("H98").Value = 4080.1304
$wsALC.Range("I98").Value = 4624.3335
$wsALC.Range("K98").Value = 4624.3335
$wsALC.Range("M98").Value = -3126.3335
$wsALC.Range("H100").Value = 14708152
$wsALC.Range("I100").Value = 21741486
$wsALC.Range("J100").Value = 2090.9092
$wsALC.Range("K100").Value = 21741486
$wsALC.Range("L100").Value = 2090.9092
$wsALC.Range("M100").Value = -21740945
$wsALC.Range("N100").Value = -3172.9092
$wsALC.Range("H112").Value = 1894.841
$wsALC.Range("I112").Value = 816.6667
$wsALC.Range("K112").Value = 2450.0001
$wsALC.Range("M112").Value = -1342.0001
$wsALC.Range("H113").Value = 2502.25
$wsALC.Range("I113").Value = 2395.1538
$wsALC.Range("J113").Value = 2966.3333
$wsALC.Range("K113").Value = 2395.1538
$wsALC.Range("L113").Value = 2966.3333
$wsALC.Range("M113").Value = 858.8462
$wsALC.Range("N113").Value = -9474.3333
$wsALC.Range("H116").Value = 2721.842
$wsALC.Range("I116").Value = 2029.2858
$wsALC.Range("J116").Value = 3125.8333
$wsALC.Range("K116").Value = 2029.2858
$wsALC.Range("L116").Value = 3125.8333
$wsALC.Range("M116").Value = 1412.7142
$wsALC.Range("N116").Value = -10009.8333
$wsALC.Range("H118").Value = 999.5
$wsALC.Range("I118").Value = 999.5
$wsALC.Range("J118").Value = 0
$wsALC.Range("K118").Value = 2998.5
$wsALC.Range("L118").Value = 0
$wsALC.Range("M118").Value = -1341.5
$wsALC.Range("N118").ClearContents()
$wsALC.Range("H122").Value = 4080.1304
$wsALC.Range("I122").Value = 4624.3335
$wsALC.Range("K122").Value = 13873.0005
$wsALC.Range("M122").Value = -11423.0005
$wsALC.Range("H132").Value = 8780712
$wsALC.Range("I132").Value = 17552288
$wsALC.Range("K132").Value = 52656864
$wsALC.Range("M132").Value = -52654334
$wsALC.Range("H135").Value = 653.5
$wsALC.Range("I135").Value = 235.1875
$wsALC.Range("J135").Value = 4000
$wsALC.Range("K135").Value = 2116.6875
$wsALC.Range("L135").Value = 36000
$wsALC.Range("M135").Value = 418.3125
$wsALC.Range("N135").Value = -41070
$wsALC.Range("H137").Value = 1000.6774
$wsALC.Range("I137").Value = 813.2909
$wsALC.Range("K137").Value = 2439.8727
$wsALC.Range("M137").Value = 110.1273000000001
$wsALC.Range("H141").Value = 556
$wsALC.Range("I141").Value = 556
$wsALC.Range("K141").Value = 1668
$wsALC.Range("M141").Value = 3512

# --- ARM ---
$wsARM.Range("H132").Value = 2332.8604
$wsARM.Range("I132").Value = 2232.2903
$wsARM.Range("J132").Value = 2592.6667
$wsARM.Range("K132").Value = 6696.8709
$wsARM.Range("L132").Value = 7778.000100000001
$wsARM.Range("M132").Value = -4166.8709
$wsARM.Range("N132").Value = -12838.0001

# --- CRP ---
$wsCRP.Range("H31").Value = 2233.5908
$wsCRP.Range("I31").Value = 2519.1177
$wsCRP.Range("J31").Value = 1262.8
$wsCRP.Range("K31").Value = 2519.1177
$wsCRP.Range("L31").Value = 1262.8
$wsCRP.Range("M31").Value = -2224.1177
$wsCRP.Range("N31").Value = -1852.8
$wsCRP.Range("H34").Value = 2233.5908
$wsCRP.Range("I34").Value = 2519.1177
$wsCRP.Range("J34").Value = 1262.8
$wsCRP.Range("K34").Value = 2519.1177
$wsCRP.Range("L34").Value = 1262.8
$wsCRP.Range("M34").Value = -2317.1177
$wsCRP.Range("N34").Value = -1666.8
$wsCRP.Range("H58").Value = 740.20514
$wsCRP.Range("I58").Value = 652.93335
$wsCRP.Range("K58").Value = 652.93335
$wsCRP.Range("M58").Value = -449.93335
$wsCRP.Range("H99").Value = 2926.5
$wsCRP.Range("I99").Value = 2889.75
$wsCRP.Range("J99").Value = 3000
$wsCRP.Range("K99").Value = 2889.75
$wsCRP.Range("L99").Value = 3000
$wsCRP.Range("M99").Value = -1391.75
$wsCRP.Range("N99").Value = -5996
$wsCRP.Range("H126").Value = 2926.5
$wsCRP.Range("I126").Value = 2889.75
$wsCRP.Range("J126").Value = 3000
$wsCRP.Range("K126").Value = 8669.25
$wsCRP.Range("L126").Value = 9000
$wsCRP.Range("M126").Value = -6199.25
$wsCRP.Range("N126").Value = -13940
$wsCRP.Range("H132").Value = 4182.5317
$wsCRP.Range("I132").Value = 4564.972
$wsCRP.Range("J132").Value = 2930.9092
$wsCRP.Range("K132").Value = 13694.916
$wsCRP.Range("L132").Value = 8792.7276
$wsCRP.Range("M132").Value = -11164.916
$wsCRP.Range("N132").Value = -13852.7276
$wsCRP.Range("H136").Value = 740.20514
$wsCRP.Range("I136").Value = 652.93335
$wsCRP.Range("K136").Value = 1958.80005
$wsCRP.Range("M136").Value = 591.1999499999999

# --- CUL ---
$wsCUL.Range("H54").Value = 1676.5
$wsCUL.Range("J54").Value = 1676.5
$wsCUL.Range("L54").Value = 5029.5
$wsCUL.Range("N54").Value = -6147.5
$wsCUL.Range("H123").Value = 2415.4583
$wsCUL.Range("I123").Value = 958.4286
$wsCUL.Range("J123").Value = 3015.4119
$wsCUL.Range("K123").Value = 2875.2858
$wsCUL.Range("L123").Value = 9046.235700000001
$wsCUL.Range("M123").Value = -425.2857999999997
$wsCUL.Range("N123").Value = -13946.2357
$wsCUL.Range("H131").Value = 16130229
$wsCUL.Range("J131").Value = 1279.9464
$wsCUL.Range("L131").Value = 3839.8392
$wsCUL.Range("N131").Value = -13919.8392

# --- GSM ---
$wsGSM.Range("H7").Value = 4296786
$wsGSM.Range("J7").Value = 2502501.5
$wsGSM.Range("L7").Value = 2502501.5
$wsGSM.Range("N7").Value = -2502725.5
$wsGSM.Range("H8").Value = 4296786
$wsGSM.Range("J8").Value = 2502501.5
$wsGSM.Range("L8").Value = 2502501.5
$wsGSM.Range("N8").Value = -2502779.5
$wsGSM.Range("H18").Value = 4000
$wsGSM.Range("I18").Value = 3000
$wsGSM.Range("J18").Value = 4500
$wsGSM.Range("K18").Value = 3000
$wsGSM.Range("L18").Value = 4500
$wsGSM.Range("M18").Value = -2707
$wsGSM.Range("N18").Value = -5086
$wsGSM.Range("H55").Value = 2940
$wsGSM.Range("I55").Value = 2000
$wsGSM.Range("J55").Value = 3880
$wsGSM.Range("K55").Value = 2000
$wsGSM.Range("L55").Value = 3880
$wsGSM.Range("M55").Value = -1673
$wsGSM.Range("N55").Value = -4534
$wsGSM.Range("H109").Value = 10213.75
$wsGSM.Range("J109").Value = 10213.75
$wsGSM.Range("L109").Value = 10213.75
$wsGSM.Range("N109").Value = -12293.75
$wsGSM.Range("H135").Value = 38831.8
$wsGSM.Range("J135").Value = 36039.75
$wsGSM.Range("L135").Value = 36039.75
$wsGSM.Range("N135").Value = -46179.75

# --- LTW ---
$wsLTW.Range("H82").Value = 1785.9656
$wsLTW.Range("I82").Value = 1685.7142
$wsLTW.Range("J82").Value = 1879.5333
$wsLTW.Range("K82").Value = 1685.7142
$wsLTW.Range("L82").Value = 1879.5333
$wsLTW.Range("M82").Value = -1324.7142
$wsLTW.Range("N82").Value = -2601.5333
$wsLTW.Range("H85").Value = 1785.9656
$wsLTW.Range("I85").Value = 1685.7142
$wsLTW.Range("J85").Value = 1879.5333
$wsLTW.Range("K85").Value = 1685.7142
$wsLTW.Range("L85").Value = 1879.5333
$wsLTW.Range("M85").Value = -437.7141999999999
$wsLTW.Range("N85").Value = -4375.5333
$wsLTW.Range("H101").Value = 15000
$wsLTW.Range("J101").Value = 15000
$wsLTW.Range("L101").Value = 15000
$wsLTW.Range("N101").Value = -21490
$wsLTW.Range("H132").Value = 19151.086
$wsLTW.Range("I132").Value = 1504.6875
$wsLTW.Range("J132").Value = 40869.73
$wsLTW.Range("K132").Value = 4514.0625
$wsLTW.Range("L132").Value = 122609.19
$wsLTW.Range("M132").Value = -1984.0625
$wsLTW.Range("N132").Value = -127669.19

# --- WVR ---
$wsWVR.Range("H132").Value = 2092.68
$wsWVR.Range("I132").Value = 1960.3405
$wsWVR.Range("K132").Value = 5881.0215
$wsWVR.Range("M132").Value = -3351.0215
